$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force specific Price (column D) cells that would otherwise be auto-parsed
# as numbers to remain stored as text, matching the source data which uses
# inline string cells for every value in columns B:E.
$textCells = @("D4", "D5", "D6", "D14", "D15", "D19", "D20", "D21", "D23", "D24", "D28", "D29", "D32", "D33", "D34", "D35", "D37", "D38", "D40", "D41", "D43", "D44", "D47", "D48", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "57.879.37"
$ws.Range("E2").Value = "  +2.26%  "
$ws.Range("D3").Value = "3.057.03"
$ws.Range("E3").Value = "  +2.25%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "526.01"
$ws.Range("E5").Value = "  +5.82%  "
$ws.Range("D6").Value = "142.81"
$ws.Range("E6").Value = "  +5.71%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +5.36%  "
$ws.Range("E9").Value = "  +6.80%  "
$ws.Range("E10").Value = "  +7.97%  "
$ws.Range("E11").Value = "  +6.00%  "
$ws.Range("E12").Value = "  +2.26%  "
$ws.Range("D13").Value = "3.579.24"
$ws.Range("E13").Value = "  +2.25%  "
$ws.Range("D14").Value = "27.11"
$ws.Range("E14").Value = "  +8.26%  "
$ws.Range("D15").Value = "0.0000169"
$ws.Range("E15").Value = "  +15.89%  "
$ws.Range("D16").Value = "57.846.66"
$ws.Range("E16").Value = "  +2.19%  "
$ws.Range("E17").Value = "  +7.98%  "
$ws.Range("D18").Value = "3.059.27"
$ws.Range("E18").Value = "  +2.37%  "
$ws.Range("D19").Value = "13.05"
$ws.Range("E19").Value = "  +5.42%  "
$ws.Range("D20").Value = "8.11"
$ws.Range("E20").Value = "  +4.49%  "
$ws.Range("D21").Value = "338.21"
$ws.Range("E21").Value = "  +4.43%  "
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").Value = "0.504"
$ws.Range("E23").Value = "  +7.70%  "
$ws.Range("D24").Value = "64.99"
$ws.Range("E24").Value = "  +5.64%  "
$ws.Range("E25").Value = "  +6.61%  "
$ws.Range("D26").Value = "0.0₃0977"
$ws.Range("E26").Value = "  +8.75%  "
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").Value = "6.93"
$ws.Range("E28").Value = "  +6.12%  "
$ws.Range("D29").Value = "7.38"
$ws.Range("E29").Value = "  +9.70%  "
$ws.Range("E30").Value = "  +5.97%  "
$ws.Range("E31").Value = "  +5.93%  "
$ws.Range("D32").Value = "21.14"
$ws.Range("E32").Value = "  +4.79%  "
$ws.Range("D33").Value = "156.47"
$ws.Range("E33").Value = "  +1.95%  "
$ws.Range("D34").Value = "4.75"
$ws.Range("E34").Value = "  +6.42%  "
$ws.Range("D35").Value = "5.99"
$ws.Range("E35").Value = "  +7.22%  "
$ws.Range("E36").Value = "  +3.87%  "
$ws.Range("D37").Value = "26.31"
$ws.Range("E37").Value = "  +14.19%  "
$ws.Range("D38").Value = "0.0702"
$ws.Range("E38").Value = "  +4.36%  "
$ws.Range("D39").Value = "3.092.46"
$ws.Range("E39").Value = "  +2.26%  "
$ws.Range("D40").Value = "37.86"
$ws.Range("E40").Value = "  +3.50%  "
$ws.Range("D41").Value = "3.89"
$ws.Range("E41").Value = "  +9.34%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "0.663"
$ws.Range("E43").Value = "  +3.67%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "1.47"
$ws.Range("E44").Value = "  +5.21%  "
$ws.Range("D45").Value = "2.325.35"
$ws.Range("E45").Value = "  +4.42%  "
$ws.Range("E46").Value = "  +4.21%  "
$ws.Range("D47").Value = "2.00"
$ws.Range("E47").Value = "  +3.03%  "
$ws.Range("D48").Value = "0.0246"
$ws.Range("E48").Value = "  +3.32%  "
$ws.Range("E49").Value = "  +5.14%  "
$ws.Range("D50").Value = "20.14"
$ws.Range("E50").Value = "  +5.84%  "
$ws.Range("D51").Value = "0.0902"
$ws.Range("E51").Value = "  +6.93%  "
